$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the first data row (row 2), pushing the
# existing data down by two rows. This makes room for the two newly
# reported IPO entries (KB스팩29호, 에이치엠씨아이비스팩7호).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Fill in the two new rows.
$ws.Range("A3").Value = "KB스팩29호"
$ws.Range("B3").Value = "2024.06.04~06.05"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
# Force the amount column to stay text (matches the rest of the column,
# which is stored as text rather than numeric) instead of letting Excel
# auto-coerce the numeric-looking string into a number.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12000"
$ws.Range("E3").ClearFormats()
$ws.Range("F3").Value = "KB증권"

$ws.Range("A4").Value = "에이치엠씨아이비스팩7호"
$ws.Range("B4").Value = "2024.06.04~06.05"
$ws.Range("C4").Value = "2,000~2,000"
$ws.Range("D4").Value = "-"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "14000"
$ws.Range("E4").ClearFormats()
$ws.Range("F4").Value = "현대차증권"

# The table keeps a rolling window of the most recent 20 IPO entries, so
# the two oldest rows (previously 코칩 / 유안타스팩16호, now pushed to rows
# 22 and 23) fall off the bottom.
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Delete()
